$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C2: filename change
$ws.Range("C2").Value = "tmpb996ygct.pdf"

# I2: AI_Резюме text update (first mentioned topic swapped)
$ws.Range("I2").Value = "Основное внимание уделяется геометрия, математическая логика. Уровень сложности: высокая сложность. Пусть 1 < p < ∞, измеримое множество E ⊂Rn, функции f, g : E →R измеримы В силу неравенства (1) |f(x..."

# M2: AI_Предварительные_знания reordered/changed
$ws.Range("M2").Value = "пределы, алгебра, тригонометрия"

# N2: AI_Области_математики reordered
$ws.Range("N2").Value = "геометрия, математическая логика, математический анализ"

# O2: AI_Рекомендации text update
$ws.Range("O2").Value = "Требует серьезной математической подготовки | Рекомендуется с преподавателем | Особое внимание уделено: геометрия, математическая логика"
